# Continued working on management. LB
$wb = $excel.ActiveWorkbook

# --- Blatt2: insert a new row at the top, shifting everything down by one ---
$ws = $wb.Worksheets.Item("Blatt2")
$ws.Rows.Item(1).Insert()

# Make room for three more rows between "PacMan" (now row 14) and "Geister"
# (now row 15) so the new "Steuerung" / "Koerner essen" / "Animation" rows fit.
$ws.Range("A15:A17").EntireRow.Insert()

# New helper cells in the freshly inserted row 1
$ws.Range("I1").Value = 48
$ws.Range("J1").Value = 49

# New content, added in the same order the author typed it so that the
# shared-strings table comes out in the expected order.
$ws.Range("B20").Value = "Levelanstieg"
$ws.Range("B21").Value = "GameOver"
$ws.Range("B22").Value = "Pause"
$ws.Range("C23").Value = "Spiel beenden"
$ws.Range("C24").Value = "Optionen"
$ws.Range("C25").Value = "Weiterspielen"
$ws.Range("B26").Value = "Zähler"
$ws.Range("B27").Value = "Leben"
$ws.Range("B28").Value = "Spielstand / Highscore"
$ws.Range("A30").Value = "Desktop-Icon"
$ws.Range("A31").Value = "Dokumentation"
$ws.Range("A32").Value = "Tests"
$ws.Range("B33").Value = "Testkonzept"
$ws.Range("B29").Value = "Hot-seat"
$ws.Range("G3").Value = "Aurel / Lukas"
$ws.Range("G4").Value = "Lukas"
$ws.Range("G12").Value = "Pascal"
$ws.Range("G13").Value = "Pascal"
$ws.Range("C15").Value = "Steuerung"
$ws.Range("C16").Value = "Körner essen"
$ws.Range("C17").Value = "Animation"

# Match the author's final selection on Blatt2
$ws.Range("A30:XFD30").Select() | Out-Null

# --- Add the new, empty "Blatt3" sheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Blatt3"
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Keep Blatt2 as the selected / active sheet, as in the original workbook.
$ws.Activate() | Out-Null
$ws.Select() | Out-Null
